# Update latest output (run 73)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update cost / unit-cost figures ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = -33.60198075
$schedule.Range("F2").Value = -0.7407844080687831
$schedule.Range("E3").Value = 372.94961925
$schedule.Range("F3").Value = 24.66598010912698

# --- Sheet "Detailed": update Price column (and two Type labels) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B7").Value = 56.97996
$detailed.Range("B8").Value = 57.06021

$detailed.Range("B9").Value = 36.2
$detailed.Range("C9").Value = "historical"

$detailed.Range("B10").Value = 56.98
$detailed.Range("C10").Value = "historical"

$detailed.Range("B11").Value = 56.98
$detailed.Range("B12").Value = 57.1

$detailed.Range("B15").Value = 50.95183
$detailed.Range("B16").Value = 24.36057
$detailed.Range("B17").Value = 0.7
$detailed.Range("B18").Value = -5.51
$detailed.Range("B19").Value = 4.14872

$detailed.Range("B21").Value = -3.6481
$detailed.Range("B22").Value = -1.5776
$detailed.Range("B23").Value = 11.7995
$detailed.Range("B24").Value = 0.00976
$detailed.Range("B25").Value = -0.93531

$detailed.Range("B27").Value = 0.00902
$detailed.Range("B28").Value = -4.49679

$detailed.Range("B30").Value = -7
$detailed.Range("B31").Value = -20.55932
$detailed.Range("B32").Value = -12.91239
$detailed.Range("B33").Value = -13.5
$detailed.Range("B34").Value = -7.10129
$detailed.Range("B35").Value = -6.44164

$detailed.Range("B37").Value = 0.66476
$detailed.Range("B38").Value = 3.98635
$detailed.Range("B39").Value = 15.83019
$detailed.Range("B40").Value = 41.35282

$detailed.Range("B42").Value = 57.3
$detailed.Range("B43").Value = 57.3

$detailed.Range("B45").Value = 46.20722
$detailed.Range("B46").Value = 56.54501
